# Commit: feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" worksheet right after "总计" (before the
#    existing "2022-Q3" sheet), populated with the Q4 fund-holding detail.
# 2. Update the "总计" (summary) sheet: insert the new 2022-Q4 row at the
#    top of the data (row 2) and shift the existing quarters down by one
#    row, adding the trailing 2021-Q3 row that is now needed at row 7.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: "总计" summary sheet — shift data down and insert 2022-Q4 row
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Row 7 does not exist yet - create it first (copy the row-above's format
# for column A so the new index cell keeps the bold/bordered style used by
# the other index cells in column A).
$summary.Range("A6").Copy()
$summary.Range("A7").PasteSpecial(-4122)

# Rewrite rows 7 down to 2 from the bottom up, shifting each quarter's
# data into the row below it, so earlier rows can be overwritten safely.
$summary.Cells.Item(7,1).Value = 5
$summary.Cells.Item(7,2).Value = "2021-Q3"
$summary.Cells.Item(7,3).Value = 11
$summary.Cells.Item(7,4).Value = 1.51

$summary.Cells.Item(6,2).Value = "2021-Q4"
$summary.Cells.Item(6,3).Value = 6
$summary.Cells.Item(6,4).Value = 0.63

$summary.Cells.Item(5,2).Value = "2022-Q1"
$summary.Cells.Item(5,3).Value = 3
$summary.Cells.Item(5,4).Value = 0.35

$summary.Cells.Item(4,2).Value = "2022-Q2"
$summary.Cells.Item(4,3).Value = 9
$summary.Cells.Item(4,4).Value = 2.48

$summary.Cells.Item(3,2).Value = "2022-Q3"
$summary.Cells.Item(3,3).Value = 19
$summary.Cells.Item(3,4).Value = 3.28

$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 10
$summary.Cells.Item(2,4).Value = 1.02

Write-Host "Updated summary sheet"

# ---------------------------------------------------------------------
# Part 2: new "2022-Q4" worksheet, inserted right before "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Match the page setup used by the rest of the workbook (values are in
# points: 0.75in/1in/0.5in).
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36
$q4.Outline.SummaryRow = 1
$q4.Outline.SummaryColumn = 1

Write-Host "Created 2022-Q4 sheet"

# Copy the header-row style (bold font + border + centered alignment,
# style index 2 in styles.xml) from the neighbouring "2022-Q3" sheet so no
# new style entries need to be introduced.
$q3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# Copy the column-A index-cell style too (also style index 2).
$q3.Range("A2").Copy()
$q4.Range("A2:A11").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Data rows: column A (index) and H (rank) are real numbers; B-G are
# stored as text (fund codes, names and percentage-looking figures are
# all text in the source data) - force text via NumberFormat, assign,
# then drop back to the default "Normal" style so no stray per-cell
# style/number-format survives on these cells.
$q4rows = @(
  @(0, "016250", "华夏远见成长一年持有混合A", "9.60", "88.62", "4.35", "0.4176", 5),
  @(1, "014062", "景顺长城专精特新量化优选股票A", "7.66", "91.15", "2.05", "0.1570", 3),
  @(2, "519975", "长信量化中小盘股票", "8.16", "93.23", "1.72", "0.1404", 5),
  @(3, "016251", "华夏远见成长一年持有混合C", "2.97", "88.62", "4.35", "0.1292", 5),
  @(4, "014063", "景顺长城专精特新量化优选股票C", "3.88", "91.15", "2.05", "0.0795", 3),
  @(5, "000458", "英大领先回报混合", "1.81", "93.66", "2.08", "0.0376", 7),
  @(6, "588160", "南方上证科创板新材料ETF", "0.90", "98.46", "3.31", "0.0298", 10),
  @(7, "588010", "博时上证科创板新材料ETF", "0.53", "98.90", "3.32", "0.0176", 10),
  @(8, "001270", "英大灵活配置混合A", "0.29", "92.68", "2.06", "0.0060", 7),
  @(9, "001271", "英大灵活配置混合B", "0.28", "92.68", "2.06", "0.0058", 7)
)

for ($i = 0; $i -lt $q4rows.Count; $i++) {
  $row = 2 + $i
  $data = $q4rows[$i]

  $q4.Cells.Item($row, 1).Value = $data[0]

  $textRange = $q4.Range("B$row`:G$row")
  $arr = New-Object 'object[,]' 1,6
  for ($c = 0; $c -lt 6; $c++) {
    $arr[0,$c] = $data[$c + 1]
  }
  $textRange.NumberFormat = "@"
  $textRange.Value = $arr
  $textRange.Style = "Normal"

  $q4.Cells.Item($row, 8).Value = $data[7]
}

Write-Host "Populated 2022-Q4 sheet"

